# Regenerate save_data: recompute column G ("K" = strikeouts) values for
# rows 2..79 to reflect the refreshed data source (K instead of Strike#).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 2
    3  = 1
    4  = 1
    5  = 3
    6  = 2
    7  = 0
    8  = 0
    9  = 0
    10 = 4
    11 = 2
    12 = 0
    13 = 1
    14 = 1
    15 = 2
    16 = 1
    17 = 2
    18 = 1
    19 = 1
    20 = 1
    21 = 2
    22 = 1
    23 = 2
    24 = 1
    25 = 1
    26 = 3
    27 = 1
    28 = 0
    29 = 1
    30 = 2
    31 = 2
    32 = 1
    33 = 1
    34 = 2
    35 = 0
    36 = 2
    37 = 1
    38 = 2
    39 = 2
    40 = 2
    41 = 3
    42 = 2
    43 = 1
    44 = 0
    45 = 1
    46 = 1
    47 = 2
    48 = 1
    49 = 0
    50 = 2
    51 = 1
    52 = 1
    53 = 1
    54 = 1
    55 = 0
    56 = 1
    57 = 2
    58 = 2
    59 = 1
    60 = 1
    61 = 3
    62 = 0
    63 = 1
    64 = 1
    65 = 1
    66 = 3
    67 = 0
    68 = 0
    69 = 2
    70 = 3
    71 = 2
    72 = 2
    73 = 2
    74 = 3
    75 = 0
    76 = 0
    77 = 0
    78 = 2
    79 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
